$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in column A (rows 3-7)
$ws.Range("A3").Value = 8
$ws.Range("A4").Value = 50
$ws.Range("A5").Value = 100
$ws.Range("A6").Value = 150
$ws.Range("A7").Value = 200

# Page setup: paper size + portrait orientation (as seen via the Page Setup dialog)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the selection to A7 (also resets any scrolled topLeftCell)
$ws.Range("A7").Select()
